# Case and Fatality Demographics Data Updated
# - Refresh the underlying counts on the three "Fatalities by ..." sheets
#   (Age Group, Gender, Race-Ethnicity) with the 11/11/21 pull.
# - Swap the tab order of "Fatalities by Age Group" and "Fatalities by Gender"
#   so Age Group now precedes Gender.
# - Leave the active selection on "Fatalities by Gender" (matching the
#   author's saved UI state) and move "Fatalities by Race-Ethnicity"'s
#   selection off of its prior cell.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update "Fatalities by Age Group" figures
# ---------------------------------------------------------------------
$wsAge = $wb.Worksheets.Item("Fatalities by Age Group")
$wsAge.Range("B3").Value = 21
$wsAge.Range("B5").Value = 631
$wsAge.Range("B6").Value = 2051
$wsAge.Range("B7").Value = 5016
$wsAge.Range("B8").Value = 9550
$wsAge.Range("B9").Value = 7265
$wsAge.Range("B10").Value = 8586
$wsAge.Range("B11").Value = 9131
$wsAge.Range("B12").Value = 8649
$wsAge.Range("B13").Value = 20406
$wsAge.Range("B15").Value = 71396

# ---------------------------------------------------------------------
# 2. Update "Fatalities by Gender" figures
# ---------------------------------------------------------------------
$wsGender = $wb.Worksheets.Item("Fatalities by Gender")
$wsGender.Range("B2").Value = 29859
$wsGender.Range("B3").Value = 41536
$wsGender.Range("B5").Value = 71396

# ---------------------------------------------------------------------
# 3. Update "Fatalities by Race-Ethnicity" figures
# ---------------------------------------------------------------------
$wsRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")
$wsRace.Range("B2").Value = 1319
$wsRace.Range("B3").Value = 7534
$wsRace.Range("B4").Value = 31202
$wsRace.Range("B5").Value = 424
$wsRace.Range("B6").Value = 30872
$wsRace.Range("B7").Value = 45

# ---------------------------------------------------------------------
# 4. Reorder tabs: "Fatalities by Age Group" moves before
#    "Fatalities by Gender" (previously Gender came first).
# ---------------------------------------------------------------------
$wsAge.Move($wsGender)

# Re-resolve sheet handles by name after the reorder, since a moved
# sheet's in-place handle now tracks the slot it vacated.
$wsAge = $wb.Worksheets.Item("Fatalities by Age Group")
$wsGender = $wb.Worksheets.Item("Fatalities by Gender")
$wsRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")

# ---------------------------------------------------------------------
# 5. Restore saved selections / active tab.
# ---------------------------------------------------------------------
$wsRace.Range("E25").Select()

$wsGender.Activate()
$wsGender.Range("C15").Select()

Write-Host "Case and Fatality Demographics Data Updated"
